# Restock inventory: set each ingredient quantity on row 2 back to full stock (1000)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:H2").Value = 1000
